# Scaling updates to better match EMEP, etc.
# - "year" sheet: change the mkd/all scaling window from 1990-2010 to
#   2000-2010 (closer to the EMEP trend) and add a dedicated Comment
#   column (shifting the old "Comment" header from H to I and giving H
#   its own select_scaling_year header).
# - add a new row for Finland (fin/all) scaling 1982-2020, skipping 1981
#   to avoid a reporting inconsistency in the inventory.

$wb = $excel.ActiveWorkbook

$wsMap    = $wb.Worksheets.Item("map")
$wsMethod = $wb.Worksheets.Item("method")
$wsYear   = $wb.Worksheets.Item("year")

# --- "year" sheet: restructure columns & update existing mkd/all row ---

# Grab the existing "Comment" header text (col H) before we overwrite it,
# so it can be moved into the new col I.
$commentHeader = $wsYear.Range("H1").Value()

$wsYear.Range("I1").Value = $commentHeader
$wsYear.Range("H1").Value = "select_scaling_year"

# mkd/all: scale from 2000 instead of 1990, leave select_scaling_year (H)
# blank (NA), and move the updated comment into the new column I.
$wsYear.Range("F2").Value = 2000
$wsYear.Range("H2").Value = "NA"

# New row: fin/all scaling 1982-2020 (skip 1981 to avoid a reporting
# inconsistency). Insert "fin" before the column-I comment text so the
# shared-string table order matches the source workbook.
$wsYear.Range("A3").Value = "fin"
$wsYear.Range("I2").Value = "Scale from 2000 so as to be closer to EMEP trend"
$wsYear.Range("B3").Value = "all"
$wsYear.Range("C3").Value = "NA"
$wsYear.Range("D3").Value = "NA"
$wsYear.Range("E3").Value = "NA"
$wsYear.Range("F3").Value = 1982
$wsYear.Range("G3").Value = 2020
$wsYear.Range("H3").Value = "NA"
$wsYear.Range("I3").Value = "Don't scale 1981 to avoid reporting inconsistency in inventory"

# Widen column F slightly now that it holds the new start-year comment
# context (cosmetic, matches the source workbook's column sizing).
$wsYear.Columns.Item(6).ColumnWidth = 11.75

# --- restore on-screen selections on each sheet (cosmetic view state) ---

$wsMap.Activate()
$excel.ActiveWindow.ScrollRow = 16
$wsMap.Range("B37").Select()

$wsMethod.Activate()
$wsMethod.Range("C35").Select()

$wsYear.Activate()
$wsYear.Rows.Item(3).Select()
